$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.196864008903503
$ws.Range("B1").Value = 2.224249839782715
$ws.Range("C1").Value = 3.518677234649658
$ws.Range("D1").Value = 2.48425817489624
$ws.Range("E1").Value = 1.128379821777344
